$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "64.070.39"
$ws.Range("E2").Value = "  +4.85%  "
$ws.Range("D3").Value = "2.755.84"
$ws.Range("E3").Value = "  +4.27%  "
$c = $ws.Range("D4")
$c.NumberFormat = "@"
$c.Value = "1.00"
$c.Style = "Normal"
$ws.Range("E4").Value = "  +0.11%  "
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = "582.21"
$c.Style = "Normal"
$ws.Range("E5").Value = "  +0.49%  "
$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = "158.52"
$c.Style = "Normal"
$ws.Range("E6").Value = "  +9.97%  "
$c = $ws.Range("D7")
$c.NumberFormat = "@"
$c.Value = "0.624"
$c.Style = "Normal"
$ws.Range("E7").Value = "  +4.22%  "
$c = $ws.Range("D8")
$c.NumberFormat = "@"
$c.Value = "1.00"
$c.Style = "Normal"
$ws.Range("E8").Value = "  +0.26%  "
$ws.Range("D9").Value = "2.752.82"
$ws.Range("E9").Value = "  +3.55%  "
$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = "6.81"
$c.Style = "Normal"
$ws.Range("E10").Value = "  +3.72%  "
$ws.Range("E11").Value = "  +3.92%  "
$c = $ws.Range("D12")
$c.NumberFormat = "@"
$c.Value = "0.393"
$c.Style = "Normal"
$ws.Range("E12").Value = "  +3.43%  "
$ws.Range("E13").Value = "  +0.90%  "
$ws.Range("D14").Value = "3.244.13"
$ws.Range("E14").Value = "  +4.48%  "
$c = $ws.Range("D15")
$c.NumberFormat = "@"
$c.Value = "27.60"
$c.Style = "Normal"
$ws.Range("E15").Value = "  +5.32%  "
$ws.Range("D16").Value = "64.003.37"
$ws.Range("E16").Value = "  +4.71%  "
$c = $ws.Range("D17")
$c.NumberFormat = "@"
$c.Value = "0.0000157"
$c.Style = "Normal"
$ws.Range("E17").Value = "  +8.39%  "
$ws.Range("D18").Value = "2.751.78"
$ws.Range("E18").Value = "  +3.82%  "
$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = "12.18"
$c.Style = "Normal"
$ws.Range("E19").Value = "  +4.60%  "
$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = "4.95"
$c.Style = "Normal"
$ws.Range("E20").Value = "  +4.31%  "
$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = "364.06"
$c.Style = "Normal"
$ws.Range("E21").Value = "  +3.07%  "
$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = "7.00"
$c.Style = "Normal"
$ws.Range("E22").Value = "  +1.73%  "
$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = "0.541"
$c.Style = "Normal"
$ws.Range("E23").Value = "  +2.29%  "
$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = "0.998"
$c.Style = "Normal"
$ws.Range("E24").Value = "  -0.06%  "
$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = "66.82"
$c.Style = "Normal"
$ws.Range("E25").Value = "  +4.14%  "
$ws.Range("E26").Value = "  +6.21%  "
$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = "8.64"
$c.Style = "Normal"
$ws.Range("E27").Value = "  +0.71%  "
$ws.Range("E28").Value = "  +0.62%  "
$ws.Range("D29").Value = "0.0₃0918"
$ws.Range("E29").Value = "  +12.38%  "
$ws.Range("E30").Value = "  +1.57%  "
$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value = "7.28"
$c.Style = "Normal"
$ws.Range("E31").Value = "  +5.43%  "
$ws.Range("E32").Value = "  +15.85%  "
$c = $ws.Range("D33")
$c.NumberFormat = "@"
$c.Value = "173.54"
$c.Style = "Normal"
$ws.Range("E33").Value = "  +2.81%  "
$c = $ws.Range("D34")
$c.NumberFormat = "@"
$c.Value = "0.999"
$c.Style = "Normal"
$ws.Range("E34").Value = "  +0.15%  "
$c = $ws.Range("D35")
$c.NumberFormat = "@"
$c.Value = "20.66"
$c.Style = "Normal"
$ws.Range("E35").Value = "  +3.05%  "
$c = $ws.Range("D36")
$c.NumberFormat = "@"
$c.Value = "4.94"
$c.Style = "Normal"
$ws.Range("E36").Value = "  +5.93%  "
$ws.Range("E37").Value = "  +6.69%  "
$ws.Range("E38").Value = "  +7.20%  "
$ws.Range("E39").Value = "  +9.45%  "
$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = "4.28"
$c.Style = "Normal"
$ws.Range("E40").Value = "  +2.95%  "
$ws.Range("E41").Value = "  -0.62%  "
$ws.Range("E42").Value = "  +15.73%  "
$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = "39.68"
$c.Style = "Normal"
$ws.Range("E43").Value = "  +3.72%  "
$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value = "22.48"
$c.Style = "Normal"
$ws.Range("E44").Value = "  +6.66%  "
$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = "21.94"
$c.Style = "Normal"
$ws.Range("E45").Value = "  +7.04%  "
$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = "0.0604"
$c.Style = "Normal"
$ws.Range("E46").Value = "  +4.80%  "
$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = "0.645"
$c.Style = "Normal"
$ws.Range("E47").Value = "  +3.37%  "
$c = $ws.Range("D48")
$c.NumberFormat = "@"
$c.Value = "0.0259"
$c.Style = "Normal"
$ws.Range("E48").Value = "  +3.79%  "
$c = $ws.Range("D49")
$c.NumberFormat = "@"
$c.Value = "137.82"
$c.Style = "Normal"
$ws.Range("E49").Value = "  +1.74%  "
$ws.Range("E50").Value = "  +2.90%  "
$ws.Range("E51").Value = "  +0.26%  "
